# Patient DB update: fix K8/L8 typing (text -> number) and append three
# new appointment rows (9, 10, 11) to the Appointments sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) K8/L8 were stored as text "12345" - convert them to real numbers.
# ---------------------------------------------------------------------
$ws.Range("K8").Value = 12345
$ws.Range("L8").Value = 12345

# ---------------------------------------------------------------------
# Helper: stamp a destination row with the same style (s="2") used by
# every other data row, by pasting the formats (only) from row 2.
# ---------------------------------------------------------------------
function Stamp-RowStyle($rowNum) {
    $ws.Range("A2:W2").Copy()
    $ws.Range("A$rowNum:W$rowNum").PasteSpecial(-4122)
}

# Give the three new rows the data-row style up front so every cell,
# including the ones we never touch (F/G), carries s="2".
Stamp-RowStyle 9
Stamp-RowStyle 10
Stamp-RowStyle 11

# ---------------------------------------------------------------------
# Row 9
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "APT_20250905_233348_001"
$ws.Range("B9").Value = "PAT_058"
$ws.Range("C9").Value = "Sachin Gupta"
$ws.Range("D9").Value = "Dr. Shreyansh"
$ws.Range("E9").Value = "Banjara Hills"
# F9 / G9 stay blank (date / time not set)
$ws.Range("H9").Value = 60
$ws.Range("I9").Value = "confirmed"
$ws.Range("J9").Value = "Blue Cross Blue Shield"
$ws.Range("K9").Value = 12345
$ws.Range("L9").Value = 12345
$ws.Range("M9").Value = "2025-09-05T23:33:48.826418"
$ws.Range("N9").Value = $false
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = "sachin@elevate.in"
$ws.Range("Q9").Value = "(701) 367-4280"
$ws.Range("R9").NumberFormat = "@"
$ws.Range("R9").Value = "12/12/1991"
$ws.Range("S9").NumberFormat = "@"
$ws.Range("S9").Value = "2025-09-08"
$ws.Range("T9").NumberFormat = "@"
$ws.Range("T9").Value = "14:00"
$ws.Range("U9").Value = "2025-09-05T23:33:48.826418"
$ws.Range("V9").Value = $true
$ws.Range("W9").Value = "2025-09-05T23:33:49.323354"

# ---------------------------------------------------------------------
# Row 10
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "APT_20250906_021436_001"
$ws.Range("B10").Value = "PAT_058"
$ws.Range("C10").Value = "Parag Patil"
$ws.Range("D10").Value = "Dr. Naresh"
$ws.Range("E10").Value = "Gachibowli"
# F10 / G10 stay blank
$ws.Range("H10").Value = 60
$ws.Range("I10").Value = "confirmed"
$ws.Range("J10").Value = "Blue Cross Blue Shield"
$ws.Range("K10").Value = 12345
$ws.Range("L10").Value = 12345
$ws.Range("M10").Value = "2025-09-06T02:14:36.455735"
$ws.Range("N10").Value = $false
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = "parag@test.com"
$ws.Range("Q10").Value = "(701) 368-4379"
$ws.Range("R10").NumberFormat = "@"
$ws.Range("R10").Value = "07/07/1990"
$ws.Range("S10").NumberFormat = "@"
$ws.Range("S10").Value = "2025-09-08"
$ws.Range("T10").NumberFormat = "@"
$ws.Range("T10").Value = "13:30"
$ws.Range("U10").Value = "2025-09-06T02:14:36.455735"
$ws.Range("V10").Value = $true
$ws.Range("W10").Value = "2025-09-06T02:14:37.023954"

# ---------------------------------------------------------------------
# Row 11
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "APT_20250906_022537_001"
$ws.Range("B11").Value = "PAT_058"
$ws.Range("C11").Value = "Parag Patil"
$ws.Range("D11").Value = "Dr. Shreyansh"
$ws.Range("E11").Value = "Jubliee Hills"
# F11 / G11 stay blank
$ws.Range("H11").Value = 60
$ws.Range("I11").Value = "confirmed"
$ws.Range("J11").Value = "Aetna"
# K11 / L11 stay text "12345" (unlike the other rows, not converted to numbers)
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = "12345"
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "12345"
$ws.Range("M11").Value = "2025-09-06T02:25:37.657630"
$ws.Range("N11").Value = $false
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = "parag@test.com"
$ws.Range("Q11").Value = "(701) 368-4370"
$ws.Range("R11").NumberFormat = "@"
$ws.Range("R11").Value = "12/12/2002"
$ws.Range("S11").NumberFormat = "@"
$ws.Range("S11").Value = "2025-09-08"
$ws.Range("T11").NumberFormat = "@"
$ws.Range("T11").Value = "14:00"
$ws.Range("U11").Value = "2025-09-06T02:25:37.657630"
$ws.Range("V11").Value = $true
$ws.Range("W11").Value = "2025-09-06T02:25:38.220367"

# ---------------------------------------------------------------------
# Re-stamp the style across the new rows: setting NumberFormat on the
# text-like date/time cells above bumped those cells onto their own
# style index. Re-pasting row 2's formats restores a uniform s="2"
# across every cell in rows 9-11 without touching any value.
# ---------------------------------------------------------------------
Stamp-RowStyle 9
Stamp-RowStyle 10
Stamp-RowStyle 11
